$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1634.8158
$ws.Range("I15").Value = 1634.8158
$ws.Range("K15").Value = 4904.4474
$ws.Range("M15").Value = -4735.4474
$ws.Range("H87").Value = 600019600
$ws.Range("J87").Value = 600019600
$ws.Range("L87").Value = 600019600
$ws.Range("N87").Value = -600022096
$ws.Range("H90").Value = 600019600
$ws.Range("J90").Value = 600019600
$ws.Range("L90").Value = 1800058800
$ws.Range("N90").Value = -1800071280
$ws.Range("H127").Value = 1696.4615
$ws.Range("I127").Value = 1198.1111
$ws.Range("J127").Value = 2817.75
$ws.Range("K127").Value = 3594.3333
$ws.Range("L127").Value = 8453.25
$ws.Range("M127").Value = 1365.6667
$ws.Range("N127").Value = -18373.25
$ws.Range("H129").Value = 1330.7
$ws.Range("I129").Value = 924.1429000000001
$ws.Range("K129").Value = 2772.4287
$ws.Range("M129").Value = 2227.5713
$ws.Range("H137").Value = 750896.9399999999
$ws.Range("I137").Value = 1055.5294
$ws.Range("K137").Value = 3166.5882
$ws.Range("M137").Value = -616.5881999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1046.069
$ws.Range("I2").Value = 416.18182
$ws.Range("K2").Value = 416.18182
$ws.Range("M2").Value = -303.18182
$ws.Range("H32").Value = 166977.52
$ws.Range("I32").Value = 166977.52
$ws.Range("K32").Value = 166977.52
$ws.Range("M32").Value = -166690.52
$ws.Range("H88").Value = 3209.6843
$ws.Range("J88").Value = 3131.5833
$ws.Range("L88").Value = 3131.5833
$ws.Range("N88").Value = -3943.5833
$ws.Range("H91").Value = 3209.6843
$ws.Range("J91").Value = 3131.5833
$ws.Range("L91").Value = 3131.5833
$ws.Range("N91").Value = -5939.5833
$ws.Range("H116").Value = 1046.069
$ws.Range("I116").Value = 416.18182
$ws.Range("K116").Value = 416.18182
$ws.Range("M116").Value = 1877.81818
$ws.Range("H122").Value = 22243528
$ws.Range("I122").Value = 22243528
$ws.Range("K122").Value = 66730584
$ws.Range("M122").Value = -66728134
$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920
$ws.Range("H132").Value = 4213.543
$ws.Range("I132").Value = 4745.759
$ws.Range("K132").Value = 14237.277
$ws.Range("M132").Value = -11707.277

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1046.069
$ws.Range("I3").Value = 416.18182
$ws.Range("K3").Value = 416.18182
$ws.Range("M3").Value = -302.18182
$ws.Range("H75").Value = 34999.5
$ws.Range("J75").Value = 49999
$ws.Range("L75").Value = 49999
$ws.Range("N75").Value = -51871
$ws.Range("H78").Value = 34999.5
$ws.Range("J78").Value = 49999
$ws.Range("L78").Value = 149997
$ws.Range("N78").Value = -159357
$ws.Range("H86").Value = 18519650
$ws.Range("I86").Value = 1120.8948
$ws.Range("J86").Value = 62501150
$ws.Range("K86").Value = 1120.8948
$ws.Range("L86").Value = 62501150
$ws.Range("M86").Value = 2.105199999999968
$ws.Range("N86").Value = -62503396
$ws.Range("H89").Value = 18519650
$ws.Range("I89").Value = 1120.8948
$ws.Range("J89").Value = 62501150
$ws.Range("K89").Value = 5604.474
$ws.Range("L89").Value = 312505750
$ws.Range("M89").Value = 11.52599999999984
$ws.Range("N89").Value = -312516982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3206
$ws.Range("I16").Value = 2429.5
$ws.Range("J16").Value = 4137.8
$ws.Range("K16").Value = 2429.5
$ws.Range("L16").Value = 4137.8
$ws.Range("M16").Value = -2142.5
$ws.Range("N16").Value = -4711.8
$ws.Range("H74").Value = 50420
$ws.Range("J74").Value = 50420
$ws.Range("L74").Value = 50420
$ws.Range("N74").Value = -52168
$ws.Range("H77").Value = 50420
$ws.Range("J77").Value = 50420
$ws.Range("L77").Value = 151260
$ws.Range("N77").Value = -159996
$ws.Range("H113").Value = 3206
$ws.Range("I113").Value = 2429.5
$ws.Range("J113").Value = 4137.8
$ws.Range("K113").Value = 2429.5
$ws.Range("L113").Value = 4137.8
$ws.Range("M113").Value = -259.5
$ws.Range("N113").Value = -8477.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2305.4583
$ws.Range("J107").Value = 1870.4706
$ws.Range("L107").Value = 5611.4118
$ws.Range("N107").Value = -9451.4118
$ws.Range("H133").Value = 7687.3125
$ws.Range("I133").Value = 6000
$ws.Range("J133").Value = 8249.75
$ws.Range("K133").Value = 18000
$ws.Range("L133").Value = 24749.25
$ws.Range("M133").Value = -12940
$ws.Range("N133").Value = -34869.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 15009.223
$ws.Range("I48").Value = 12762.45
$ws.Range("K48").Value = 12762.45
$ws.Range("M48").Value = -12277.45
$ws.Range("H49").Value = 29296.6
$ws.Range("J49").Value = 29296.6
$ws.Range("L49").Value = 29296.6
$ws.Range("N49").Value = -29664.6
$ws.Range("H70").Value = 23812522
$ws.Range("I70").Value = 30305964
$ws.Range("J70").Value = 3239
$ws.Range("K70").Value = 30305964
$ws.Range("L70").Value = 3239
$ws.Range("M70").Value = -30305694
$ws.Range("N70").Value = -3779
$ws.Range("H73").Value = 23812522
$ws.Range("I73").Value = 30305964
$ws.Range("J73").Value = 3239
$ws.Range("K73").Value = 30305964
$ws.Range("L73").Value = 3239
$ws.Range("M73").Value = -30305028
$ws.Range("N73").Value = -5111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10999.5
$ws.Range("I3").Value = 9999
$ws.Range("J3").Value = 12000
$ws.Range("K3").Value = 9999
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = -9887
$ws.Range("N3").Value = -12224
$ws.Range("H15").Value = 10999.5
$ws.Range("I15").Value = 9999
$ws.Range("J15").Value = 12000
$ws.Range("K15").Value = 9999
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = -9829
$ws.Range("N15").Value = -12340
$ws.Range("H46").Value = 1410.8889
$ws.Range("I46").Value = 785.5238000000001
$ws.Range("J46").Value = 3599.6667
$ws.Range("K46").Value = 785.5238000000001
$ws.Range("L46").Value = 3599.6667
$ws.Range("M46").Value = -597.5238000000001
$ws.Range("N46").Value = -3975.6667
$ws.Range("H132").Value = 3426.7837
$ws.Range("I132").Value = 2531.7576
$ws.Range("J132").Value = 10810.75
$ws.Range("K132").Value = 7595.2728
$ws.Range("L132").Value = 32432.25
$ws.Range("M132").Value = -5065.2728
$ws.Range("N132").Value = -37492.25
$ws.Range("H136").Value = 24631.066
$ws.Range("I136").Value = 3323.3044
$ws.Range("K136").Value = 9969.913199999999
$ws.Range("M136").Value = -7419.913199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4139.75
$ws.Range("I96").Value = 1503
$ws.Range("K96").Value = 1503
$ws.Range("M96").Value = -130
